$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H32").Value = 2591.6
$ws.Range("I32").Value = 986
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 986
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -660
$ws.Range("N32").Value = -5652
$ws.Range("H116").Value = 4461.5
$ws.Range("I116").Value = 3333
$ws.Range("J116").Value = 5590
$ws.Range("K116").Value = 3333
$ws.Range("L116").Value = 5590
$ws.Range("M116").Value = 109
$ws.Range("N116").Value = -12474
$ws.Range("H132").Value = 31315.285
$ws.Range("I132").Value = 31315.285
$ws.Range("K132").Value = 93945.855
$ws.Range("M132").Value = -91415.855
$ws.Range("H137").Value = 3404.0527
$ws.Range("I137").Value = 2268.2
$ws.Range("J137").Value = 4666.1113
$ws.Range("K137").Value = 6804.599999999999
$ws.Range("L137").Value = 13998.3339
$ws.Range("M137").Value = -4254.599999999999
$ws.Range("N137").Value = -19098.3339
$ws.Range("H140").Value = 50780
$ws.Range("J140").Value = 50780
$ws.Range("L140").Value = 50780
$ws.Range("N140").Value = -61140
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 942.2353000000001
$ws.Range("I2").Value = 809.7692
$ws.Range("J2").Value = 1372.75
$ws.Range("K2").Value = 809.7692
$ws.Range("L2").Value = 1372.75
$ws.Range("M2").Value = -696.7692
$ws.Range("N2").Value = -1598.75
$ws.Range("H61").Value = 3561.5789
$ws.Range("I61").Value = 2282.3845
$ws.Range("J61").Value = 6333.1665
$ws.Range("K61").Value = 2282.3845
$ws.Range("L61").Value = 6333.1665
$ws.Range("M61").Value = -2070.3845
$ws.Range("N61").Value = -6757.1665
$ws.Range("H102").Value = 6875
$ws.Range("I102").Value = 3750
$ws.Range("K102").Value = 3750
$ws.Range("M102").Value = -2128
$ws.Range("H110").Value = 500.36365
$ws.Range("I110").Value = 480.44446
$ws.Range("K110").Value = 480.44446
$ws.Range("M110").Value = 1564.55554
$ws.Range("H116").Value = 942.2353000000001
$ws.Range("I116").Value = 809.7692
$ws.Range("J116").Value = 1372.75
$ws.Range("K116").Value = 809.7692
$ws.Range("L116").Value = 1372.75
$ws.Range("M116").Value = 1484.2308
$ws.Range("N116").Value = -5960.75
$ws.Range("H131").Value = 69994.5
$ws.Range("J131").Value = 69994.5
$ws.Range("L131").Value = 69994.5
$ws.Range("N131").Value = -80074.5
$ws.Range("H136").Value = 3561.5789
$ws.Range("I136").Value = 2282.3845
$ws.Range("J136").Value = 6333.1665
$ws.Range("K136").Value = 6847.1535
$ws.Range("L136").Value = 18999.4995
$ws.Range("M136").Value = -4297.1535
$ws.Range("N136").Value = -24099.4995
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 942.2353000000001
$ws.Range("I3").Value = 809.7692
$ws.Range("J3").Value = 1372.75
$ws.Range("K3").Value = 809.7692
$ws.Range("L3").Value = 1372.75
$ws.Range("M3").Value = -695.7692
$ws.Range("N3").Value = -1600.75
$ws.Range("H36").Value = 2137
$ws.Range("I36").Value = 2137
$ws.Range("K36").Value = 2137
$ws.Range("M36").Value = -1603
$ws.Range("H39").Value = 9999
$ws.Range("J39").Value = 9999
$ws.Range("L39").Value = 9999
$ws.Range("N39").Value = -10777
$ws.Range("H80").Value = 796.5
$ws.Range("I80").Value = 598.625
$ws.Range("J80").Value = 1060.3334
$ws.Range("K80").Value = 598.625
$ws.Range("L80").Value = 1060.3334
$ws.Range("M80").Value = 399.375
$ws.Range("N80").Value = -3056.3334
$ws.Range("H83").Value = 796.5
$ws.Range("I83").Value = 598.625
$ws.Range("J83").Value = 1060.3334
$ws.Range("K83").Value = 2993.125
$ws.Range("L83").Value = 5301.666999999999
$ws.Range("M83").Value = 1998.875
$ws.Range("N83").Value = -15285.667
$ws.Range("H134").Value = 2162.6667
$ws.Range("I134").Value = 2162.6667
$ws.Range("K134").Value = 6488.000100000001
$ws.Range("M134").Value = -3953.000100000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 46664.668
$ws.Range("J20").Value = 46664.668
$ws.Range("L20").Value = 46664.668
$ws.Range("N20").Value = -47136.668
$ws.Range("H30").Value = 46664.668
$ws.Range("J30").Value = 46664.668
$ws.Range("L30").Value = 46664.668
$ws.Range("N30").Value = -46846.668
$ws.Range("H31").Value = 6723.357
$ws.Range("I31").Value = 997
$ws.Range("K31").Value = 997
$ws.Range("M31").Value = -702
$ws.Range("H34").Value = 6723.357
$ws.Range("I34").Value = 997
$ws.Range("K34").Value = 997
$ws.Range("M34").Value = -795
$ws.Range("H122").Value = 3098.25
$ws.Range("I122").Value = 1196.5
$ws.Range("K122").Value = 3589.5
$ws.Range("M122").Value = -1139.5
$ws.Range("H128").Value = 46664.668
$ws.Range("J128").Value = 46664.668
$ws.Range("L128").Value = 46664.668
$ws.Range("N128").Value = -56624.668
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 540.4167
$ws.Range("I14").Value = 540.4167
$ws.Range("K14").Value = 1621.2501
$ws.Range("M14").Value = -1448.2501
$ws.Range("H37").Value = 111920.75
$ws.Range("J37").Value = 111920.75
$ws.Range("L37").Value = 335762.25
$ws.Range("N37").Value = -335986.25
$ws.Range("H69").Value = 2837
$ws.Range("I69").Value = 3255.5
$ws.Range("K69").Value = 9766.5
$ws.Range("M69").Value = -8955.5
$ws.Range("H72").Value = 2837
$ws.Range("I72").Value = 3255.5
$ws.Range("K72").Value = 29299.5
$ws.Range("M72").Value = -25243.5
$ws.Range("H132").Value = 1706.5
$ws.Range("J132").Value = 1623.3
$ws.Range("L132").Value = 14609.7
$ws.Range("N132").Value = -19669.7
$ws.Range("H133").Value = 130
$ws.Range("I133").Value = 130
$ws.Range("K133").Value = 390
$ws.Range("M133").Value = 4670
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 49999
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 49999
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 49999
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -50991
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 7423.3076
$ws.Range("I100").Value = 5500
$ws.Range("J100").Value = 8278.111000000001
$ws.Range("K100").Value = 5500
$ws.Range("L100").Value = 8278.111000000001
$ws.Range("M100").Value = -4959
$ws.Range("N100").Value = -9360.111000000001
$ws.Range("H136").Value = 3771.25
$ws.Range("I136").Value = 2292.5
$ws.Range("K136").Value = 6877.5
$ws.Range("M136").Value = -4327.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 35756.92
$ws.Range("I4").Value = 34570
$ws.Range("K4").Value = 34570
$ws.Range("M4").Value = -34457
$ws.Range("H122").Value = 4997.5
$ws.Range("J122").Value = 4995
$ws.Range("L122").Value = 14985
$ws.Range("N122").Value = -19885
